$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The three SVM classification-report rows (4-6) were re-run and their
# order/ratings refreshed: SVM polynomial, SVM linear, SVM radial (was
# SVM radial, SVM polynomial, SVM linear), with the "0.0" rating format
# now carried by the top row instead of the second one.
$ws.Range("A4").Value = "SVM polynomial"
$ws.Range("B4").Value = 85
$ws.Range("B4").NumberFormat = "0.0"

$ws.Range("A5").Value = "SVM linear"
$ws.Range("B5").Value = 85.1
$ws.Range("B5").ClearFormats()

$ws.Range("A6").Value = "SVM radial"
$ws.Range("B6").Value = 85.2

# The rating (value) axis on the comparison chart now shows one decimal
# place instead of the default General format.
$chart = $ws.ChartObjects(1).Chart
$chart.Axes(2).TickLabels.NumberFormat = "0.0"

# Reflect the cell that was last being worked on when the file was saved.
$ws.Range("C4").Select() | Out-Null
